$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue 'D2' '37.766.18'
$ws.Range('E2').Value = '  +0.64%  '
Set-TextValue 'D3' '2.115.52'
$ws.Range('E3').Value = '  +1.87%  '
$ws.Range('E4').Value = '  -0.06%  '
Set-TextValue 'D5' '236.11'
$ws.Range('E5').Value = '  +0.34%  '
Set-TextValue 'D6' '0.628'
$ws.Range('E6').Value = '  +1.01%  '
Set-TextValue 'D7' '58.43'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  -0.08%  '
Set-TextValue 'D9' '0.392'
$ws.Range('E9').Value = '  +1.63%  '
Set-TextValue 'D10' '0.0785'
$ws.Range('E10').Value = '  +2.64%  '
$ws.Range('E11').Value = '  +1.07%  '
Set-TextValue 'D12' '2.426.41'
$ws.Range('E12').Value = '  +1.77%  '
Set-TextValue 'D13' '14.59'
$ws.Range('E13').Value = '  +0.47%  '
Set-TextValue 'D14' '21.31'
$ws.Range('E14').Value = '  +1.06%  '
Set-TextValue 'D15' '0.790'
$ws.Range('E15').Value = '  +1.04%  '
Set-TextValue 'D16' '5.24'
$ws.Range('E16').Value = '  +0.55%  '
Set-TextValue 'D17' '2.104.62'
$ws.Range('E17').Value = '  +1.29%  '
Set-TextValue 'D18' '37.730.07'
$ws.Range('E18').Value = '  +0.08%  '
Set-TextValue 'D19' '6.22'
$ws.Range('E19').Value = '  +0.08%  '
Set-TextValue 'D20' '70.30'
$ws.Range('E20').Value = '  +1.06%  '
Set-TextValue 'D21' '0.0₃0826'
$ws.Range('E21').Value = '  +1.06%  '
Set-TextValue 'D22' '227.61'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  +0.07%  '
Set-TextValue 'D24' '2.42'
$ws.Range('E24').Value = '  -0.81%  '
Set-TextValue 'D25' '2.43'
$ws.Range('E25').Value = '  -3.58%  '
Set-TextValue 'D26' '168.54'
$ws.Range('E26').Value = '  +0.86%  '
Set-TextValue 'D27' '9.00'
$ws.Range('E27').Value = '  -0.52%  '
Set-TextValue 'D28' '0.136'
$ws.Range('E28').Value = '  +4.20%  '
Set-TextValue 'D29' '1.43'
$ws.Range('E29').Value = '  -4.54%  '
Set-TextValue 'D30' '19.48'
$ws.Range('E30').Value = '  +1.09%  '
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('E32').Value = '  +2.68%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D33' '2.59'
$ws.Range('E33').Value = '  -0.14%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D34' '0.0623'
$ws.Range('E34').Value = '  -0.60%  '
Set-TextValue 'D35' '4.59'
$ws.Range('E35').Value = '  -0.28%  '
Set-TextValue 'D36' '3.54'
$ws.Range('E36').Value = '  +5.57%  '
$ws.Range('E37').Value = '  +1.31%  '
Set-TextValue 'D38' '0.999'
$ws.Range('E38').Value = '  -0.12%  '
Set-TextValue 'D39' '5.67'
$ws.Range('E39').Value = '  -5.60%  '
$ws.Range('E40').Value = '  -0.22%  '
Set-TextValue 'D41' '0.0966'
$ws.Range('E41').Value = '  +1.41%  '
Set-TextValue 'D42' '97.93'
$ws.Range('E42').Value = '  +2.43%  '
Set-TextValue 'D43' '1.479.68'
$ws.Range('E43').Value = '  +0.65%  '
Set-TextValue 'D44' '0.0214'
$ws.Range('E44').Value = '  +0.69%  '
$ws.Range('E45').Value = '  -0.85%  '
Set-TextValue 'D46' '4.23'
$ws.Range('E46').Value = '  -6.88%  '
Set-TextValue 'D47' '1.06'
$ws.Range('E47').Value = '  +1.69%  '
$ws.Range('E48').Value = '  -1.22%  '
Set-TextValue 'D49' '7.34'
$ws.Range('E49').Value = '  +1.32%  '
$ws.Range('E50').Value = '  +3.07%  '
Set-TextValue 'D51' '2.312.03'
$ws.Range('E51').Value = '  +1.82%  '
